# ELG4913 BOM.xlsx -- "updated with made purchases"
#
# The sheet is reorganised from a single price list into two sections:
#   - rows 3-12  : "still to buy" list, with a Total / Per person subtotal
#   - rows 16-21 : "Ordered" list (parts that were actually purchased,
#                  generally at a slightly different price than planned),
#                  with its own Total: / Per person: subtotal.
#
# We rebuild the sheet content/formula/style state directly rather than
# trying to replay literal cut/paste gestures, since the end state is
# what matters.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$currencyFmt = "$#,##0.00;[Red]-$#,##0.00"

$ALIEXPRESS = "https://www.aliexpress.com/item/1005005275736468.html?spm=a2g0o.productlist.main.1.3a1072cewm2hWM&algo_pvid=b268c8e1-af98-4201-8858-09cd4e6887e0&algo_exp_id=b268c8e1-af98-4201-8858-09cd4e6887e0-0&pdp_npi=4%40dis%21CAD%2173.14%2171.68%21%21%2152.61%2151.56%21%40210307bf17265994460883461ec49b%2112000038407879035%21sea%21CA%210%21ABX&curPageLogUid=dF82d1Pl71AG&utparam-url=scene%3Asearch%7Cquery_from%3A"
$AMAZON_ETHERNET_NEW = "https://www.amazon.ca/Ethernet-Higher-Bandwidth-Internet-Network/dp/B017P34W6C/ref=sr_1_5?crid=2UEA4J2GFAALW&dib=eyJ2IjoiMSJ9.XcDYd-kLRvKFIX-X1VJTtW3GFvx8he748TcX6fsPJmWtuwR9pf-4hz9_LugBUMfHQ4QrpkUJhKPcq9Oh8XAd6Mfv5uS-zGj7f-Dn5G9lnedaomx2WezBxwMFXgEDohA9hXs9rWHY5sP6rwzSGwTXVdMopK-SpBzaNGm2LsnhGS45dK1TBsd6kWi4pqm86bC3lj05sSNQaOYEapFXo6lR3hlSxWYQ4fGUobp05n2GcA2n4y4QaVhLrkj4Vyydlptk3NovA-UjLubbjSF-rolckk7xZ-DQMyNpDZbWPprUHqQ.KOOQIVYrd1Bd43431-r4OP3yD9vFSjODHJY8-IOfeWk&dib_tag=se&keywords=ethernet+cable+100ft&qid=1726690261&sprefix=ethernet+cable+100ft%2Caps%2C105&sr=8-5"
$AMAZON_B0BXWJFCVJ = "https://www.amazon.ca/gp/product/B0BXWJFCVJ/ref=ox_sc_act_title_2?smid=A1XEC9TMFJSNSW&psc=1"
$MOUSER = "https://www.mouser.ca/ProductDetail/Pimoroni/PIM366?qs=lc2O%252BfHJPVbvcWNSB8Ff5Q%3D%3D"

# --- start clean: drop the two existing live hyperlinks (Aliexpress /
#     Pishop) -- only the Aliexpress one survives, re-pointed at the
#     "Ordered" section further down the sheet.
$ws.Hyperlinks.Delete()

# --- clear out the old row 13 / row 14 entries entirely (their content
#     is being relocated into the new "Ordered" block) so no stale rows
#     remain below the rebuilt "to buy" list.
$ws.Range("C13:E14").Clear()

# ======================== "still to buy" section ========================

# Header row is unchanged.
$ws.Range("B2").Value = "ELG4913 BOM"
$ws.Range("C3").Value = "Component"
$ws.Range("D3").Value = "Price"
$ws.Range("E3").Value = "Purchase Link"

# Row 4: Ethernet cable, re-priced, pointing at the new Amazon listing
# (plain text, not a live hyperlink).
$ws.Range("C4").Value = "Ethernet Cable(100ft)"
$ws.Range("D4").Value = 23.99
$ws.Range("D4").NumberFormat = $currencyFmt
$ws.Range("E4").Value = $AMAZON_ETHERNET_NEW

# Row 5: Suspension, still TBD.
$ws.Range("C5").Value = "Suspension"
$ws.Range("D5").Value = "TBD"

# Row 6/7: Total / Per person subtotal for the "to buy" list (bold,
# right-aligned -- new style).
$ws.Range("C6").Value = "Total"
$ws.Range("C6").Font.Bold = $true
$ws.Range("C6").HorizontalAlignment = -4152
$ws.Range("D6").Formula = "=D4"
$ws.Range("D6").NumberFormat = $currencyFmt

$ws.Range("C7").Value = "Per person"
$ws.Range("C7").Font.Bold = $true
$ws.Range("C7").HorizontalAlignment = -4152
$ws.Range("D7").Formula = "=D6/5"
$ws.Range("D7").NumberFormat = $currencyFmt

# Row 8: blank spacer row, carrying the same bold/right style + currency
# format as rows 6/7/20/21 but no value.
$ws.Range("C8").Font.Bold = $true
$ws.Range("C8").HorizontalAlignment = -4152
$ws.Range("D8").NumberFormat = $currencyFmt

# Row 9: "Things to ask about:" header, moved down from row 11.
$ws.Range("C9").Value = "Things to ask about:"
$ws.Range("C9").Font.Bold = $true
$ws.Range("D9").Value = 1.15
$ws.Range("D9").NumberFormat = $currencyFmt

# Rows 10-12: remaining "to buy" parts, moved down from rows 12-14.
$ws.Range("C10").Value = "ON/OFF Switch"
$ws.Range("D10").Value = 20.32
$ws.Range("D10").NumberFormat = $currencyFmt

$ws.Range("C11").Value = "ADS1115 ADC 16-bit ADC"
$ws.Range("D11").Value = 2.69
$ws.Range("D11").NumberFormat = $currencyFmt

$ws.Range("C12").Value = "Voltage Regulator 7.4V-5V"

# ============================ "Ordered" section ============================

$ws.Range("C16").Value = "Ordered"
$ws.Range("C16").Font.Bold = $true

# Row 17: Geophone, actually purchased via the Aliexpress link (keeps the
# live hyperlink + the matching "display" text).
$ws.Range("C17").Value = "Geophone SM-24"
$ws.Range("D17").Value = 87.78
$ws.Range("D17").NumberFormat = $currencyFmt
$ws.Range("E17").Value = $ALIEXPRESS
$ws.Hyperlinks.Add($ws.Range("E17"), $ALIEXPRESS, "", "", $ALIEXPRESS) | Out-Null

# Row 18: ADS1115, purchased from Amazon (plain text link, no live
# hyperlink).
$ws.Range("C18").Value = "ADS1115 ADC 16-bit ADC"
$ws.Range("D18").Value = 14.69
$ws.Range("D18").NumberFormat = $currencyFmt
$ws.Range("E18").Value = $AMAZON_B0BXWJFCVJ

# Row 19: MLX90640 camera, purchased from Mouser instead of Pishop.
# Keeps the wrap-text style from its old row 5 home, two-line tall, and
# the "looks like a hyperlink" blue/underlined style on the link cell,
# even though it is no longer a real clickable hyperlink.
$ws.Range("C19").Value = "MLX90640 Thermal Camera Breakout - Wide Angle"
$ws.Range("C19").WrapText = $true
$ws.Range("C19").RowHeight = 30
$ws.Range("D19").Value = 113.88
$ws.Range("D19").NumberFormat = $currencyFmt
$ws.Range("E19").Value = $MOUSER
$ws.Range("E19").Style = "Hyperlink"

# Row 20/21: Total: / Per person: subtotal for the "Ordered" list (same
# bold/right style as the "to buy" subtotal above).
$ws.Range("C20").Value = "Total:"
$ws.Range("C20").Font.Bold = $true
$ws.Range("C20").HorizontalAlignment = -4152
$ws.Range("D20").Formula = "=D17+D18+D19"
$ws.Range("D20").NumberFormat = $currencyFmt

$ws.Range("C21").Value = "Per person:"
$ws.Range("C21").Font.Bold = $true
$ws.Range("C21").HorizontalAlignment = -4152
$ws.Range("D21").Formula = "=D20/5"
$ws.Range("D21").NumberFormat = $currencyFmt

# Matches the saved selection in the target workbook.
$ws.Range("C12").Select() | Out-Null
